$d = $word.ActiveDocument

$oldPath = "Program: /home/rstudio/work/cctu/tests/testthat/analysis.R"
$newPath = "Program: /home/rstudio/Documents/GitHub/cctu/tests/testthat/analysis.R"

foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute(
                $oldPath, $true, $false, $false, $false, $false,
                $true, 1, $false, $newPath, 2)
        }
    }
}
